# Auto-generated PowerShell COM-interop script
# Applies the BRVM daily data refresh described by the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet 'Recommandations' ----
$wsReco = $wb.Worksheets.Item("Recommandations")

$recoData = @(
    @('BRVM - SERVICES PUBLICS', 0, 8, 3383.27, 111.89, '🟡 Observer', '➖ Neutre'),
    @('AIR LIQUIDE CI', 0, 4, 2850.0, 715.0, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI', 0, 4, 2840.0, 700.0, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AUTRES SECTEURS', 0, 4, 2466.62, 618.56, '🟡 Observer', '➖ Neutre'),
    @('BRVM - DISTRIBUTION', 0, 4, 2024.83, 522.18, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AGRICULTURE', 0, 4, 1483.31, 365.17, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TRANSPORT', 0, 4, 1440.49, 363.49, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 719.46, 187.95, '🟡 Observer', '➖ Neutre'),
    @('BRVM - FINANCES', 0, 4, 560.77, 140.19, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 4, 560.25, 140.62, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 4, 551.12, 137.78, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 4, 515.59, 129.32, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 4, 436.22, 108.82, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL                    (**)', 0, 2, 408.63, 204.3, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 4, 384.43, 97.09, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE                 (**)', 0, 1, 219.23, 219.23, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE              (**)', 0, 1, 219.1, 219.1, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE                (**)', 0, 1, 214.18, 214.18, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL                 (**)', 0, 1, 204.98, 204.98, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE        (**)', 0, 1, 194.82, 194.82, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE           (**)', 0, 1, 193.79, 193.79, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE         (**)', 0, 1, 189.95, 189.95, '🟡 Observer', '➖ Neutre'),
    @('TRACTAFRIC MOTORS CI (PRSC)', 3, 0, 9.43, 4.08, '🟢 Achat', '✅ Renforcer'),
    @('CIE CI (CIEC)', 2, 0, 7.43, 2.89, '🟡 Observer', '➖ Neutre'),
    @('SERVAIR ABIDJAN CI (ABJC)', 2, 1, 5.93, -5.17, '🟡 Observer', '👀 À surveiller'),
    @('UNIWAX CI (UNXC)', 1, 0, 5.88, 5.88, '🟡 Observer', '➖ Neutre'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 1, 0, 5.0, 5.0, '🟡 Observer', '➖ Neutre'),
    @('ORANGE COTE D''IVOIRE (ORAC)', 1, 0, 3.47, 3.47, '🟡 Observer', '➖ Neutre'),
    @('SOGB CI (SOGC)', 1, 0, 2.07, 2.07, '🟡 Observer', '➖ Neutre'),
    @('SICOR CI (SICC)', 1, 2, 1.82, -2.86, '🟡 Observer', '👀 À surveiller'),
    @('CFAO MOTORS CI (CFAC)', 2, 2, 1.57, 7.26, '🟡 Observer', '👀 À surveiller'),
    @('NEI-CEDA CI (NEIC)', 1, 0, 1.45, 1.45, '🟡 Observer', '➖ Neutre'),
    @('SAPH CI (SPHC)', 1, 0, 1.44, 1.44, '🟡 Observer', '➖ Neutre'),
    @('SETAO CI (STAC)', 1, 1, 0.19, 7.27, '🟡 Observer', '👀 À surveiller'),
    @('TOTAL', 0, 4, 0.0, 0.0, '🟡 Observer', '➖ Neutre'),
    @('ONATEL BF (ONTBF)', 1, 1, -0.65, 5.49, '🟡 Observer', '👀 À surveiller'),
    @('PALM CI (PALC)', 0, 1, -1.99, -1.99, '🟡 Observer', '➖ Neutre'),
    @('VIVO ENERGY CI (SHEC)', 0, 1, -1.99, -1.99, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA NG (BOAN)', 0, 1, -2.5, -2.5, '🟡 Observer', '➖ Neutre'),
    @('TOTALENERGIES MARKETING CI (TTLC)', 1, 1, -3.0, 4.35, '🟡 Observer', '👀 À surveiller'),
    @('ECOBANK COTE D''''IVOIRE (ECOC)', 0, 1, -3.85, -3.85, '🟡 Observer', '➖ Neutre'),
    @('BERNABE CI (BNBC)', 0, 1, -4.46, -4.46, '🟡 Observer', '➖ Neutre'),
    @('SAFCA CI (SAFC)', 1, 2, -7.37, 7.44, '🟡 Observer', '👀 À surveiller'),
    @('SOLIBRA CI (SLBC)', 0, 1, -7.46, -7.46, '🟡 Observer', '➖ Neutre'),
    @('AIR LIQUIDE CI (SIVC)', 0, 2, -9.79, -6.29, '🟡 Observer', '➖ Neutre'),
    @('SICABLE CI (CABC)', 0, 2, -10.62, -3.25, '🟡 Observer', '➖ Neutre')
)

$r = 2
foreach ($row in $recoData) {
    $wsReco.Cells.Item($r, 1).Value = $row[0]
    $wsReco.Cells.Item($r, 2).Value = $row[1]
    $wsReco.Cells.Item($r, 3).Value = $row[2]
    $wsReco.Cells.Item($r, 4).Value = $row[3]
    $wsReco.Cells.Item($r, 5).Value = $row[4]
    $wsReco.Cells.Item($r, 6).Value = $row[5]
    $wsReco.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# The dataset now has one fewer row (48 -> 47 data rows); clear the
# now-unused former last row (old dimension was A1:G48).
$wsReco.Range("A48:G48").Clear()

# ---- Sheet 'Top_YTD' ----
$wsYtd = $wb.Worksheets.Item("Top_YTD")

$ytdData = @(
    @('BRVM - SERVICES PUBLICS', 9288624.0),
    @('AIR LIQUIDE CI', 435631.6),
    @('NEI-CEDA CI', 430220.0),
    @('BRVM - AUTRES SECTEURS', 263662.25),
    @('BRVM - DISTRIBUTION', 134837.78),
    @('BRVM - AGRICULTURE', 49030.58),
    @('BRVM - TRANSPORT', 44715.17),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 6030.04),
    @('BRVM - FINANCES', 3228.39),
    @('BRVM-PRESTIGE', 3221.16)
)

$r = 2
foreach ($row in $ytdData) {
    $wsYtd.Cells.Item($r, 1).Value = $row[0]
    $wsYtd.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

Write-Host "BRVM recommandations + Top_YTD refreshed"
